$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.8546833333333335
$ws.Range("H2").Value = 2.56405
$ws.Range("I2").Value = 0.3097546281380014
$ws.Range("J2").Value = 0.3097546281380015
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 6.949630179233334
$ws.Range("R2").Value = 62.5466716131
$ws.Range("S2").Value = 0.006476955561577232
$ws.Range("T2").Value = 0.006476955561577235
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.8546833333333335
$ws.Range("H3").Value = 2.56405
$ws.Range("I3").Value = 0.3097546281380014
$ws.Range("J3").Value = 0.3097546281380015
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 208.0097516670889
$ws.Range("R3").Value = 1872.087765003801
$ws.Range("S3").Value = 0.1938621024681745
$ws.Range("T3").Value = 0.1938621024681746
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.8546833333333335
$ws.Range("H4").Value = 2.56405
$ws.Range("I4").Value = 0.3097546281380014
$ws.Range("J4").Value = 0.3097546281380015
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 88.83804060619445
$ws.Range("R4").Value = 799.5423654557502
$ws.Range("S4").Value = 0.08279577852981408
$ws.Range("T4").Value = 0.08279577852981411
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.8546833333333335
$ws.Range("H5").Value = 2.56405
$ws.Range("I5").Value = 0.3097546281380014
$ws.Range("J5").Value = 0.3097546281380015
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 28.56244807604445
$ws.Range("R5").Value = 257.0620326844
$ws.Range("S5").Value = 0.02661979157843555
$ws.Range("T5").Value = 0.02661979157843556
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.251276
$ws.Range("H6").Value = 3.753828
$ws.Range("I6").Value = 0.453487879032787
$ws.Range("J6").Value = 0.4534878790327871
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 10.174417954584
$ws.Range("R6").Value = 91.569761591256
$ws.Range("S6").Value = 0.009482411474738925
$ws.Range("T6").Value = 0.009482411474738927
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.251276
$ws.Range("H7").Value = 3.753828
$ws.Range("I7").Value = 0.453487879032787
$ws.Range("J7").Value = 0.4534878790327871
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 304.5310466180321
$ws.Range("R7").Value = 2740.779419562289
$ws.Range("S7").Value = 0.2838185637502789
$ws.Range("T7").Value = 0.283818563750279
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.251276
$ws.Range("H8").Value = 3.753828
$ws.Range("I8").Value = 0.453487879032787
$ws.Range("J8").Value = 0.4534878790327871
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 130.06092872318
$ws.Range("R8").Value = 1170.54835850862
$ws.Range("S8").Value = 0.1212149184793647
$ws.Range("T8").Value = 0.1212149184793647
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.251276
$ws.Range("H9").Value = 3.753828
$ws.Range("I9").Value = 0.453487879032787
$ws.Range("J9").Value = 0.4534878790327871
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 41.816078990816
$ws.Range("R9").Value = 376.3447109173441
$ws.Range("S9").Value = 0.0389719853284045
$ws.Range("T9").Value = 0.03897198532840451
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1689346666666667
$ws.Range("H10").Value = 0.506804
$ws.Range("I10").Value = 0.06122536009783416
$ws.Range("J10").Value = 0.06122536009783416
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 1.373647305378667
$ws.Range("R10").Value = 12.362825748408
$ws.Range("S10").Value = 0.001280219569208708
$ws.Range("T10").Value = 0.001280219569208708
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1689346666666667
$ws.Range("H11").Value = 0.506804
$ws.Range("I11").Value = 0.06122536009783416
$ws.Range("J11").Value = 0.06122536009783416
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 41.11471078328712
$ws.Range("R11").Value = 370.0323970495841
$ws.Range("S11").Value = 0.03831832022748415
$ws.Range("T11").Value = 0.03831832022748415
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1689346666666667
$ws.Range("H12").Value = 0.506804
$ws.Range("I12").Value = 0.06122536009783416
$ws.Range("J12").Value = 0.06122536009783416
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 17.55951495929556
$ws.Range("R12").Value = 158.03563463366
$ws.Range("S12").Value = 0.01636521586631458
$ws.Range("T12").Value = 0.01636521586631458
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1689346666666667
$ws.Range("H13").Value = 0.506804
$ws.Range("I13").Value = 0.06122536009783416
$ws.Range("J13").Value = 0.06122536009783416
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 5.645585279043556
$ws.Range("R13").Value = 50.81026751139201
$ws.Range("S13").Value = 0.005261604434826718
$ws.Range("T13").Value = 0.005261604434826719
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.4843330000000001
$ws.Range("H14").Value = 1.452999
$ws.Range("I14").Value = 0.1755321327313773
$ws.Range("J14").Value = 0.1755321327313773
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 3.938224956922
$ws.Range("R14").Value = 35.444024612298
$ws.Range("S14").Value = 0.003670369124633357
$ws.Range("T14").Value = 0.003670369124633357
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.4843330000000001
$ws.Range("H15").Value = 1.452999
$ws.Range("I15").Value = 0.1755321327313773
$ws.Range("J15").Value = 0.1755321327313773
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 117.8752212954227
$ws.Range("R15").Value = 1060.876991658804
$ws.Range("S15").Value = 0.1098580140887093
$ws.Range("T15").Value = 0.1098580140887093
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.4843330000000001
$ws.Range("H16").Value = 1.452999
$ws.Range("I16").Value = 0.1755321327313773
$ws.Range("J16").Value = 0.1755321327313773
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 50.34284985189834
$ws.Range("R16").Value = 453.0856486670851
$ws.Range("S16").Value = 0.0469188133648101
$ws.Range("T16").Value = 0.04691881336481011
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.4843330000000001
$ws.Range("H17").Value = 1.452999
$ws.Range("I17").Value = 0.1755321327313773
$ws.Range("J17").Value = 0.1755321327313773
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 16.18580312086133
$ws.Range("R17").Value = 145.672228087752
$ws.Range("S17").Value = 0.0150849361532245
$ws.Range("T17").Value = 0.0150849361532245
